# Fix "Recorded By" (column G) attribution strings.
#
# Several rows in column G contain a comma-separated list of recorder
# names/emails such as "dnasr281@gmail.com, System" or
# "system, backup@backdoor.com, System". The list order needs to be
# reversed wherever the word "System"/"system" appears in the list
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# Lists that do not mention "system" at all (e.g.
# "admin@admin.com, dnasr281@gmail.com") are left untouched, as are
# single-value cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }
    if ($val.ToLower().IndexOf("system") -lt 0) { continue }

    $parts = $val.Split(",")
    $count = $parts.Length

    $rebuilt = ""
    for ($i = $count - 1; $i -ge 0; $i--) {
        $piece = $parts[$i].Trim()
        if ($i -eq ($count - 1)) {
            $rebuilt = $piece
        } else {
            $rebuilt = $rebuilt + ", " + $piece
        }
    }

    $cell.Value = $rebuilt
    $changed++
}

Write-Host "Recorded By entries reordered:" $changed
